$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1177974704067688, 1.452459299700094, 8.919389290004966, 2.986534662448264, 3.051279910712037, 23)
    3  = @(0.02651855621502405, 1.404872730675466, 8.509181565190792, 2.917050147870412, 2.98557263622811, 22)
    4  = @(-0.4315838055284564, 0.8785582649148802, 3.741529435360165, 1.934303346261947, 1.932104577106019, 21)
    5  = @(-0.02783490465184608, 0.6090290681569026, 0.9013455341055104, 0.9493921919341397, 0.9736371053067127, 20)
    6  = @(-0.01094518282552908, 0.6512406853883913, 0.945541919671934, 0.9723897982146532, 0.9989722588254877, 19)
    7  = @(-0.1131319738049607, 0.4987054024040059, 0.4823568942290646, 0.6945191820454383, 0.7051092912023117, 18)
    8  = @(-0.02891905992389743, 0.4566464661289574, 0.38608079962584, 0.6213540050774921, 0.6397829851394089, 17)
    9  = @(0.06681065338993029, 0.4080515307318274, 0.312996434393734, 0.5594608425919851, 0.5736737944045387, 16)
    10 = @(0.0231535205448071, 0.3284333828178208, 0.2480728887976249, 0.4980691606570566, 0.5149932072547927, 15)
    11 = @(0.02456212780780078, 0.3234927182903493, 0.2262431411680219, 0.4756502298622612, 0.4929470093722863, 14)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
